$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.333.54"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "1.880.95"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4828"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2894"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").Value = "1.880.49"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.182"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6609"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "30.303.62"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007755"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.435"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "2.134.49"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "199.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.177"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.425"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.435"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.265"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.043"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05047"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7416"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.149"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01842"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.632"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9139"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.901"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4320"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.648"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1351"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.563"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -14.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.895"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05716"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.98%  "
